$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: update title only (link stays the same)
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 36: update title and link
$ws.Range("D36").Value = "Applications of Self-Supervised Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/355"

# Row 46: update title and link
$ws.Range("D46").Value = "[한국생명공학연구원] 2022년 03월, 생물정보학(Bioinformatics 채용), 바이오데이터 수집/관리/분석 분야 정규직"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/431"
